$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.021072
$ws.Range("H2").Value = 0.042144
$ws.Range("I2").Value = 0.006826597554061716
$ws.Range("J2").Value = 0.004661306488643927
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 0.0006892229760000001
$ws.Range("R2").Value = 0.002756891904
$ws.Range("S2").Value = 0.006826597554061716
$ws.Range("T2").Value = 0.004661306488643927

$ws.Range("H3").Value = 0.7921320000000001
$ws.Range("I3").Value = 0.08554110310196809
$ws.Range("J3").Value = 0.08761318411784576
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 0.008636351152
$ws.Range("R3").Value = 0.05181810691200001
$ws.Range("S3").Value = 0.08554110310196809
$ws.Range("T3").Value = 0.08761318411784576

$ws.Range("G4").Value = 0.08179599999999999
$ws.Range("H4").Value = 0.245388
$ws.Range("I4").Value = 0.02649906859965984
$ws.Range("J4").Value = 0.02714096138561494
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 0.002675383568
$ws.Range("R4").Value = 0.016052301408
$ws.Range("S4").Value = 0.02649906859965984
$ws.Range("T4").Value = 0.02714096138561494

$ws.Range("G5").Value = 0.197936
$ws.Range("H5").Value = 0.395872
$ws.Range("I5").Value = 0.06412440268891229
$ws.Range("J5").Value = 0.04378513482992712
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.006474090688
$ws.Range("R5").Value = 0.025896362752
$ws.Range("S5").Value = 0.06412440268891229
$ws.Range("T5").Value = 0.04378513482992712

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.439847666666667
$ws.Range("H6").Value = 7.319543
$ws.Range("I6").Value = 0.7904260684106801
$ws.Range("J6").Value = 0.8095727334806435
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 0.07980253748133334
$ws.Range("R6").Value = 0.478815224888
$ws.Range("S6").Value = 0.7904260684106801
$ws.Range("T6").Value = 0.8095727334806435

$ws.Range("G7").Value = 0.08205433333333334
$ws.Range("H7").Value = 0.246163
$ws.Range("I7").Value = 0.02658275964471802
$ws.Range("J7").Value = 0.02722667969732477
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 0.002683833134666667
$ws.Range("R7").Value = 0.016102998808
$ws.Range("S7").Value = 0.02658275964471802
$ws.Range("T7").Value = 0.02722667969732477
